# "fixed PER bug": the Team/PER sheet's Team column (B) was showing the
# wrong team per row, and the PER column (C) held stale, unscaled values.
# Rewrite both columns for rows 2-30 with the corrected team labels and
# recomputed per-game PER figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "POR"
$ws.Cells.Item(2, 3).Value = 13.66666666666666
$ws.Cells.Item(3, 2).Value = "NJN"
$ws.Cells.Item(3, 3).Value = 9.926666666666668
$ws.Cells.Item(4, 2).Value = "CLE"
$ws.Cells.Item(4, 3).Value = 12.41666666666667
$ws.Cells.Item(5, 2).Value = "DAL"
$ws.Cells.Item(5, 3).Value = 13.4375
$ws.Cells.Item(6, 2).Value = "MIA"
$ws.Cells.Item(6, 3).Value = 11.49285714285714
$ws.Cells.Item(7, 2).Value = "SEA"
$ws.Cells.Item(7, 3).Value = 10.91
$ws.Cells.Item(8, 2).Value = "ATL"
$ws.Cells.Item(8, 3).Value = 8.822222222222221
$ws.Cells.Item(9, 2).Value = "WAS"
$ws.Cells.Item(9, 3).Value = 13.76666666666667
$ws.Cells.Item(10, 2).Value = "MIL"
$ws.Cells.Item(10, 3).Value = 13.08181818181818
$ws.Cells.Item(11, 2).Value = "LAC"
$ws.Cells.Item(11, 3).Value = 13.34666666666666
$ws.Cells.Item(12, 2).Value = "SAS"
$ws.Cells.Item(12, 3).Value = 12.62307692307692
$ws.Cells.Item(13, 2).Value = "DET"
$ws.Cells.Item(13, 3).Value = 11.93333333333333
$ws.Cells.Item(14, 2).Value = "ORL"
$ws.Cells.Item(14, 3).Value = 13.575
$ws.Cells.Item(15, 2).Value = "UTA"
$ws.Cells.Item(15, 3).Value = 13.21538461538462
$ws.Cells.Item(16, 2).Value = "MEM"
$ws.Cells.Item(16, 3).Value = 12.90769230769231
$ws.Cells.Item(17, 2).Value = "HOU"
$ws.Cells.Item(17, 3).Value = 13.32307692307692
$ws.Cells.Item(18, 2).Value = "DEN"
$ws.Cells.Item(18, 3).Value = 11.61538461538462
$ws.Cells.Item(19, 2).Value = "LAL"
$ws.Cells.Item(19, 3).Value = 12.56428571428571
$ws.Cells.Item(20, 2).Value = "GSW"
$ws.Cells.Item(20, 3).Value = 12.21875
$ws.Cells.Item(21, 2).Value = "IND"
$ws.Cells.Item(21, 3).Value = 13.91333333333333
$ws.Cells.Item(22, 2).Value = "CHI"
$ws.Cells.Item(22, 3).Value = 12.57333333333333
$ws.Cells.Item(23, 2).Value = "PHI"
$ws.Cells.Item(23, 3).Value = 14.33636363636364
$ws.Cells.Item(24, 2).Value = "BOS"
$ws.Cells.Item(24, 3).Value = 10.825
$ws.Cells.Item(25, 2).Value = "TOR"
$ws.Cells.Item(25, 3).Value = 12.94117647058824
$ws.Cells.Item(26, 2).Value = "SAC"
$ws.Cells.Item(26, 3).Value = 13.72142857142857
$ws.Cells.Item(27, 2).Value = "PHO"
$ws.Cells.Item(27, 3).Value = 11.47142857142857
$ws.Cells.Item(28, 2).Value = "NOH"
$ws.Cells.Item(28, 3).Value = 12.67857142857143
$ws.Cells.Item(29, 2).Value = "NYK"
$ws.Cells.Item(29, 3).Value = 12.66153846153846
$ws.Cells.Item(30, 2).Value = "MIN"
$ws.Cells.Item(30, 3).Value = 14.74615384615385